$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.480.25"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "2.935.97"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'590.46"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").Value = "'147.62"
$ws.Range("E6").Value = "  +6.20%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.506"
$ws.Range("E8").Value = "  +2.84%  "
$ws.Range("D9").Value = "2.938.34"
$ws.Range("E9").Value = "  +1.87%  "
$ws.Range("D10").Value = "'7.12"
$ws.Range("E10").Value = "  +2.76%  "
$ws.Range("E11").Value = "  +9.15%  "
$ws.Range("E12").Value = "  +2.16%  "
$ws.Range("E13").Value = "  +7.83%  "
$ws.Range("D14").Value = "'32.29"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").Value = "3.425.15"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("D17").Value = "62.484.87"
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'6.63"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "2.939.88"
$ws.Range("E19").Value = "  +2.91%  "
$ws.Range("D20").Value = "'434.65"
$ws.Range("E20").Value = "  +2.38%  "
$ws.Range("D21").Value = "'13.44"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").Value = "'0.662"
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").Value = "'80.20"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").Value = "'11.12"
$ws.Range("E25").Value = "  +7.24%  "
$ws.Range("D26").Value = "'11.88"
$ws.Range("E26").Value = "  +4.86%  "
$ws.Range("D27").Value = "'2.10"
$ws.Range("E27").Value = "  +2.61%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "'7.26"
$ws.Range("E29").Value = "  +9.80%  "
$ws.Range("E30").Value = "  +4.62%  "
$ws.Range("D31").Value = "'2.58"
$ws.Range("E31").Value = "  +1.97%  "
$ws.Range("E32").Value = "  +21.41%  "
$ws.Range("E33").Value = "  +6.12%  "
$ws.Range("E34").Value = "  +1.79%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'0.990"
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("D38").Value = "'3.05"
$ws.Range("E38").Value = "  +8.77%  "
$ws.Range("D39").Value = "'49.58"
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("D40").Value = "'2.01"
$ws.Range("E40").Value = "  +5.64%  "
$ws.Range("D41").Value = "'0.116"
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("E42").Value = "  +0.44%  "
$ws.Range("D43").Value = "'0.276"
$ws.Range("E43").Value = "  +4.33%  "
$ws.Range("D44").Value = "'39.63"
$ws.Range("E44").Value = "  +5.53%  "
$ws.Range("D45").Value = "2.697.22"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("D46").Value = "'135.37"
$ws.Range("E46").Value = "  +3.08%  "
$ws.Range("E47").Value = "  +3.39%  "
$ws.Range("D48").Value = "'354.06"
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("D51").Value = "'22.54"
$ws.Range("E51").Value = "  +1.46%  "
